$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect temporarily so the values can be updated.
$ws.Unprotect()

# Update the confidential disclosure date string (2021-05-20 -> 2021-05-21)
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-21 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.2460174370354169
$ws.Range("E2").Value = -0.002512983749371811

$ws.Range("D3").Value = 0.4997960832406509
$ws.Range("E3").Value = 0.002900843881856519

$ws.Range("D4").Value = 0.09641878005471372
$ws.Range("E4").Value = -0.0007911392405063333

$ws.Range("D5").Value = 0.1015799131004597
$ws.Range("E5").Value = 0.0008203445447088065

$ws.Range("D6").Value = 0.05618778656875888
$ws.Range("E6").Value = 0.002974147792267345

$ws.Range("E7").Value = 0.001005753217429239

# Restore sheet protection to match the original workbook's protected state.
$ws.Protect()
